$d = $word.ActiveDocument
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$d.Content.Find.Execute("2023-04-27 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-04-28 Friday", 2) | Out-Null
$d.Content.Find.Execute("21+27=", $true, $false, $false, $false, $false, $true, 1, $false, "75+20=", 2) | Out-Null
$d.Content.Find.Execute("6+20=", $true, $false, $false, $false, $false, $true, 1, $false, "96-92=", 2) | Out-Null
$d.Content.Find.Execute("85+4=", $true, $false, $false, $false, $false, $true, 1, $false, "79-47=", 2) | Out-Null
$d.Content.Find.Execute("77-41=", $true, $false, $false, $false, $false, $true, 1, $false, "70-1=", 2) | Out-Null
$d.Content.Find.Execute("35+18=", $true, $false, $false, $false, $false, $true, 1, $false, "31-11=", 2) | Out-Null
$d.Content.Find.Execute("16+66=", $true, $false, $false, $false, $false, $true, 1, $false, "38-13=", 2) | Out-Null
$d.Content.Find.Execute("2+89=", $true, $false, $false, $false, $false, $true, 1, $false, "88-64=", 2) | Out-Null
$d.Content.Find.Execute("51+14=", $true, $false, $false, $false, $false, $true, 1, $false, "83-62=", 2) | Out-Null
$d.Content.Find.Execute("78+6=", $true, $false, $false, $false, $false, $true, 1, $false, "92-73=", 2) | Out-Null
$d.Content.Find.Execute("52+0=", $true, $false, $false, $false, $false, $true, 1, $false, "47-9=", 2) | Out-Null
$d.Content.Find.Execute("16+79=", $true, $false, $false, $false, $false, $true, 1, $false, "14-5=", 2) | Out-Null
$d.Content.Find.Execute("92-16=", $true, $false, $false, $false, $false, $true, 1, $false, "18+36=", 2) | Out-Null
$d.Content.Find.Execute("95-17=", $true, $false, $false, $false, $false, $true, 1, $false, "44+12=", 2) | Out-Null
$d.Content.Find.Execute("45-41=", $true, $false, $false, $false, $false, $true, 1, $false, "71-54=", 2) | Out-Null
$d.Content.Find.Execute("41+2=", $true, $false, $false, $false, $false, $true, 1, $false, "59-37=", 2) | Out-Null
$d.Content.Find.Execute("10+45=", $true, $false, $false, $false, $false, $true, 1, $false, "64-51=", 2) | Out-Null
$d.Content.Find.Execute("17-4=", $true, $false, $false, $false, $false, $true, 1, $false, "33+20=", 2) | Out-Null
$d.Content.Find.Execute("32-27=", $true, $false, $false, $false, $false, $true, 1, $false, "99-54=", 2) | Out-Null
$d.Content.Find.Execute("0+94=", $true, $false, $false, $false, $false, $true, 1, $false, "31+32=", 2) | Out-Null
$d.Content.Find.Execute("28-12=", $true, $false, $false, $false, $false, $true, 1, $false, "9+37=", 2) | Out-Null
$d.Content.Find.Execute("82-47=", $true, $false, $false, $false, $false, $true, 1, $false, "3+39=", 2) | Out-Null
$d.Content.Find.Execute("4+79=", $true, $false, $false, $false, $false, $true, 1, $false, "22+4=", 2) | Out-Null
$d.Content.Find.Execute("43+56=", $true, $false, $false, $false, $false, $true, 1, $false, "56+5=", 2) | Out-Null
$d.Content.Find.Execute("57-20=", $true, $false, $false, $false, $false, $true, 1, $false, "71-27=", 2) | Out-Null
$d.Content.Find.Execute("12+35=", $true, $false, $false, $false, $false, $true, 1, $false, "62-7=", 2) | Out-Null
$d.Content.Find.Execute("10+18=", $true, $false, $false, $false, $false, $true, 1, $false, "54+44=", 2) | Out-Null
$d.Content.Find.Execute("14+52=", $true, $false, $false, $false, $false, $true, 1, $false, "40-29=", 2) | Out-Null
$d.Content.Find.Execute("58+38=", $true, $false, $false, $false, $false, $true, 1, $false, "47+37=", 2) | Out-Null
$d.Content.Find.Execute("87-23=", $true, $false, $false, $false, $false, $true, 1, $false, "34+40=", 2) | Out-Null
$d.Content.Find.Execute("79-23=", $true, $false, $false, $false, $false, $true, 1, $false, "62-16=", 2) | Out-Null
$d.Content.Find.Execute("19+32=", $true, $false, $false, $false, $false, $true, 1, $false, "5+81=", 2) | Out-Null
$d.Content.Find.Execute("8+55=", $true, $false, $false, $false, $false, $true, 1, $false, "43-18=", 2) | Out-Null
$d.Content.Find.Execute("83-79=", $true, $false, $false, $false, $false, $true, 1, $false, "71+24=", 2) | Out-Null
$d.Content.Find.Execute("54+45=", $true, $false, $false, $false, $false, $true, 1, $false, "87-73=", 2) | Out-Null
$d.Content.Find.Execute("4+22=", $true, $false, $false, $false, $false, $true, 1, $false, "11-10=", 2) | Out-Null
$d.Content.Find.Execute("71-31=", $true, $false, $false, $false, $false, $true, 1, $false, "16+65=", 2) | Out-Null
$d.Content.Find.Execute("23+0=", $true, $false, $false, $false, $false, $true, 1, $false, "65+2=", 2) | Out-Null
$d.Content.Find.Execute("5+44=", $true, $false, $false, $false, $false, $true, 1, $false, "5+6=", 2) | Out-Null
$d.Content.Find.Execute("62+7=", $true, $false, $false, $false, $false, $true, 1, $false, "92-52=", 2) | Out-Null
$d.Content.Find.Execute("26+49=", $true, $false, $false, $false, $false, $true, 1, $false, "41-2=", 2) | Out-Null
$d.Content.Find.Execute("31+67=", $true, $false, $false, $false, $false, $true, 1, $false, "31+0=", 2) | Out-Null
$d.Content.Find.Execute("80+16=", $true, $false, $false, $false, $false, $true, 1, $false, "52+39=", 2) | Out-Null
$d.Content.Find.Execute("37+57=", $true, $false, $false, $false, $false, $true, 1, $false, "38-34=", 2) | Out-Null
$d.Content.Find.Execute("58-43=", $true, $false, $false, $false, $false, $true, 1, $false, "82-32=", 2) | Out-Null
$d.Content.Find.Execute("33+42=", $true, $false, $false, $false, $false, $true, 1, $false, "9+6=", 2) | Out-Null
$d.Content.Find.Execute("35+60=", $true, $false, $false, $false, $false, $true, 1, $false, "98-4=", 2) | Out-Null
$d.Content.Find.Execute("24-3=", $true, $false, $false, $false, $false, $true, 1, $false, "74-39=", 2) | Out-Null
$d.Content.Find.Execute("63-5=", $true, $false, $false, $false, $false, $true, 1, $false, "24-20=", 2) | Out-Null
$d.Content.Find.Execute("84-71=", $true, $false, $false, $false, $false, $true, 1, $false, "18+75=", 2) | Out-Null
$d.Content.Find.Execute("57-1=", $true, $false, $false, $false, $false, $true, 1, $false, "15+71=", 2) | Out-Null
$d.Content.Find.Execute("87+5=", $true, $false, $false, $false, $false, $true, 1, $false, "79-42=", 2) | Out-Null
$d.Content.Find.Execute("91-41=", $true, $false, $false, $false, $false, $true, 1, $false, "24+8=", 2) | Out-Null
$d.Content.Find.Execute("7+85=", $true, $false, $false, $false, $false, $true, 1, $false, "12+38=", 2) | Out-Null
$d.Content.Find.Execute("58+21=", $true, $false, $false, $false, $false, $true, 1, $false, "37+20=", 2) | Out-Null
$d.Content.Find.Execute("96-88=", $true, $false, $false, $false, $false, $true, 1, $false, "20+7=", 2) | Out-Null
$d.Content.Find.Execute("81-52=", $true, $false, $false, $false, $false, $true, 1, $false, "1+1=", 2) | Out-Null
$d.Content.Find.Execute("31-20=", $true, $false, $false, $false, $false, $true, 1, $false, "74-14=", 2) | Out-Null
$d.Content.Find.Execute("26-24=", $true, $false, $false, $false, $false, $true, 1, $false, "80-66=", 2) | Out-Null
$d.Content.Find.Execute("48-36=", $true, $false, $false, $false, $false, $true, 1, $false, "2+68=", 2) | Out-Null
$d.Content.Find.Execute("11-7=", $true, $false, $false, $false, $false, $true, 1, $false, "20+16=", 2) | Out-Null
$d.Content.Find.Execute("6+57=", $true, $false, $false, $false, $false, $true, 1, $false, "37-19=", 2) | Out-Null
$d.Content.Find.Execute("8+39=", $true, $false, $false, $false, $false, $true, 1, $false, "50+27=", 2) | Out-Null
$d.Content.Find.Execute("59+15=", $true, $false, $false, $false, $false, $true, 1, $false, "50+28=", 2) | Out-Null
$d.Content.Find.Execute("10+77=", $true, $false, $false, $false, $false, $true, 1, $false, "20-5=", 2) | Out-Null
$d.Content.Find.Execute("60-28=", $true, $false, $false, $false, $false, $true, 1, $false, "59-49=", 2) | Out-Null
$d.Content.Find.Execute("86-66=", $true, $false, $false, $false, $false, $true, 1, $false, "91-46=", 2) | Out-Null
$d.Content.Find.Execute("94-33=", $true, $false, $false, $false, $false, $true, 1, $false, "8+30=", 2) | Out-Null
$d.Content.Find.Execute("12+79=", $true, $false, $false, $false, $false, $true, 1, $false, "21+15=", 2) | Out-Null
$d.Content.Find.Execute("51+18=", $true, $false, $false, $false, $false, $true, 1, $false, "68+15=", 2) | Out-Null
$d.Content.Find.Execute("94-86=", $true, $false, $false, $false, $false, $true, 1, $false, "94-38=", 2) | Out-Null
$d.Content.Find.Execute("21-17=", $true, $false, $false, $false, $false, $true, 1, $false, "94-62=", 2) | Out-Null
$d.Content.Find.Execute("28-11=", $true, $false, $false, $false, $false, $true, 1, $false, "14+64=", 2) | Out-Null
$d.Content.Find.Execute("37+0=", $true, $false, $false, $false, $false, $true, 1, $false, "7+39=", 2) | Out-Null
$d.Content.Find.Execute("1+8=", $true, $false, $false, $false, $false, $true, 1, $false, "6+43=", 2) | Out-Null
$d.Content.Find.Execute("61+25=", $true, $false, $false, $false, $false, $true, 1, $false, "54-28=", 2) | Out-Null
$d.Content.Find.Execute("84-61=", $true, $false, $false, $false, $false, $true, 1, $false, "62-60=", 2) | Out-Null
$d.Content.Find.Execute("39+5=", $true, $false, $false, $false, $false, $true, 1, $false, "72+2=", 2) | Out-Null
$d.Content.Find.Execute("31-18=", $true, $false, $false, $false, $false, $true, 1, $false, "95-30=", 2) | Out-Null
$d.Content.Find.Execute("36+14=", $true, $false, $false, $false, $false, $true, 1, $false, "9+2=", 2) | Out-Null
$d.Content.Find.Execute("8+63=", $true, $false, $false, $false, $false, $true, 1, $false, "50+16=", 2) | Out-Null
$d.Content.Find.Execute("42-14=", $true, $false, $false, $false, $false, $true, 1, $false, "4+17=", 2) | Out-Null
$d.Content.Find.Execute("73-27=", $true, $false, $false, $false, $false, $true, 1, $false, "35+32=", 2) | Out-Null
$d.Content.Find.Execute("90-0=", $true, $false, $false, $false, $false, $true, 1, $false, "50-17=", 2) | Out-Null
$d.Content.Find.Execute("28+56=", $true, $false, $false, $false, $false, $true, 1, $false, "30+34=", 2) | Out-Null
$d.Content.Find.Execute("57-49=", $true, $false, $false, $false, $false, $true, 1, $false, "39+36=", 2) | Out-Null
$d.Content.Find.Execute("40+1=", $true, $false, $false, $false, $false, $true, 1, $false, "33+20=", 2) | Out-Null
$d.Content.Find.Execute("70+15=", $true, $false, $false, $false, $false, $true, 1, $false, "40-25=", 2) | Out-Null
$d.Content.Find.Execute("31-8=", $true, $false, $false, $false, $false, $true, 1, $false, "72+23=", 2) | Out-Null
$d.Content.Find.Execute("76-46=", $true, $false, $false, $false, $false, $true, 1, $false, "76-32=", 2) | Out-Null
$d.Content.Find.Execute("48-1=", $true, $false, $false, $false, $false, $true, 1, $false, "41+14=", 2) | Out-Null
$d.Content.Find.Execute("40+15=", $true, $false, $false, $false, $false, $true, 1, $false, "49-42=", 2) | Out-Null
$d.Content.Find.Execute("14+75=", $true, $false, $false, $false, $false, $true, 1, $false, "85-47=", 2) | Out-Null
$d.Content.Find.Execute("60-9=", $true, $false, $false, $false, $false, $true, 1, $false, "5+7=", 2) | Out-Null
$d.Content.Find.Execute("8+69=", $true, $false, $false, $false, $false, $true, 1, $false, "6+83=", 2) | Out-Null
$d.Content.Find.Execute("53+46=", $true, $false, $false, $false, $false, $true, 1, $false, "23+8=", 2) | Out-Null
$d.Content.Find.Execute("48-29=", $true, $false, $false, $false, $false, $true, 1, $false, "36-17=", 2) | Out-Null
$d.Content.Find.Execute("32+60=", $true, $false, $false, $false, $false, $true, 1, $false, "19-17=", 2) | Out-Null
$d.Content.Find.Execute("40+58=", $true, $false, $false, $false, $false, $true, 1, $false, "82-62=", 2) | Out-Null
$d.Content.Find.Execute("27+15=", $true, $false, $false, $false, $false, $true, 1, $false, "95-19=", 2) | Out-Null
$d.Content.Find.Execute("46+2=", $true, $false, $false, $false, $false, $true, 1, $false, "54+3=", 2) | Out-Null
